$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($CellRef, $TextValue, [bool]$Protect = $false)
    $rng = $ws.Range($CellRef)
    if ($Protect) {
        $origStyle = $rng.Style
        $rng.Value = "'" + $TextValue
        $rng.Style = $origStyle
    } else {
        $rng.Value = $TextValue
    }
}

Set-CellText "D2" "26.147.79" $false
Set-CellText "E2" "  +0.48%  " $false
Set-CellText "D3" "1.655.75" $false
Set-CellText "E3" "  +0.21%  " $false
Set-CellText "E4" "  +0.11%  " $false
Set-CellText "D5" "218.07" $true
Set-CellText "E5" "  +0.66%  " $false
Set-CellText "D6" "0.5306" $true
Set-CellText "E6" "  +1.83%  " $false
Set-CellText "E7" "  +0.16%  " $false
Set-CellText "D8" "0.2615" $true
Set-CellText "E8" "  +0.02%  " $false
Set-CellText "E9" "  +1.07%  " $false
Set-CellText "D10" "20.43" $true
Set-CellText "E10" "  -0.76%  " $false
Set-CellText "D11" "0.07808" $true
Set-CellText "E11" "  +0.92%  " $false
Set-CellText "B12" "Polkadot" $false
Set-CellText "C12" "https://coinranking.com/coin/25W7FG7om+polkadot-dot" $false
Set-CellText "D12" "4.517" $true
Set-CellText "E12" "  +1.26%  " $false
Set-CellText "B13" "WrappedEther" $false
Set-CellText "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" $false
Set-CellText "D13" "1.623.50" $false
Set-CellText "E13" "  -1.30%  " $false
Set-CellText "D14" "1.883.36" $false
Set-CellText "E14" "  +0.26%  " $false
Set-CellText "D15" "0.5490" $true
Set-CellText "E15" "  +1.30%  " $false
Set-CellText "D16" "0.0₅8206" $false
Set-CellText "E16" "  +1.36%  " $false
Set-CellText "D17" "65.40" $true
Set-CellText "D18" "26.131.90" $false
Set-CellText "E18" "  +0.41%  " $false
Set-CellText "D20" "4.605" $true
Set-CellText "E20" "  +0.86%  " $false
Set-CellText "D21" "191.46" $true
Set-CellText "E21" "  +0.08%  " $false
Set-CellText "E22" "  +0.78%  " $false
Set-CellText "D23" "6.021" $true
Set-CellText "E23" "  +0.76%  " $false
Set-CellText "D24" "1.005" $true
Set-CellText "E24" "  +0.13%  " $false
Set-CellText "D25" "145.19" $true
Set-CellText "E25" "  +5.08%  " $false
Set-CellText "D27" "7.216" $true
Set-CellText "E27" "  -0.39%  " $false
Set-CellText "D28" "15.98" $true
Set-CellText "E28" "  -0.79%  " $false
Set-CellText "D29" "1.466" $true
Set-CellText "E29" "  +5.03%  " $false
Set-CellText "D30" "0.05764" $true
Set-CellText "E30" "  -3.27%  " $false
Set-CellText "E31" "  +0.18%  " $false
Set-CellText "D32" "3.562" $true
Set-CellText "E32" "  +1.75%  " $false
Set-CellText "D33" "3.271" $true
Set-CellText "E33" "  +0.98%  " $false
Set-CellText "D34" "1.600" $true
Set-CellText "E34" "  +2.58%  " $false
Set-CellText "D35" "2.801" $true
Set-CellText "E35" "  +1.71%  " $false
Set-CellText "D36" "0.9509" $true
Set-CellText "E36" "  +0.29%  " $false
Set-CellText "E37" "  +0.18%  " $false
Set-CellText "D38" "0.5745" $true
Set-CellText "E38" "  +1.16%  " $false
Set-CellText "D39" "0.01610" $true
Set-CellText "E39" "  +0.85%  " $false
Set-CellText "D40" "0.8536" $true
Set-CellText "E40" "  +0.99%  " $false
Set-CellText "D41" "5.808" $true
Set-CellText "E41" "  -1.30%  " $false
Set-CellText "D42" "104.70" $true
Set-CellText "E42" "  +3.92%  " $false
Set-CellText "B43" "PaxDollar" $false
Set-CellText "C43" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp" $false
Set-CellText "D43" "1.004" $true
Set-CellText "E43" "  +0.24%  " $false
Set-CellText "B44" "Maker" $false
Set-CellText "C44" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr" $false
Set-CellText "D44" "1.042.11" $false
Set-CellText "E44" "  +3.75%  " $false
Set-CellText "D45" "1.797.61" $false
Set-CellText "D46" "56.90" $true
Set-CellText "E46" "  +0.50%  " $false
Set-CellText "D47" "1.005" $true
Set-CellText "E47" "  +0.69%  " $false
Set-CellText "E48" "  +1.03%  " $false
Set-CellText "D49" "7.849" $true
Set-CellText "E49" "  -1.07%  " $false
Set-CellText "D50" "0.05151" $true
Set-CellText "E50" "  +0.04%  " $false
Set-CellText "D51" "1.442" $true
Set-CellText "E51" "  -2.07%  " $false
